$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert two new blank rows at the very top. Everything (data, the Séance
#    table, the progress-legend fragments, rich text, drawings) shifts down
#    by two rows. This alone reproduces the target layout for every row
#    from (old) row 6 onward landing at (new) row 8 onward.
# ---------------------------------------------------------------------------
$ws.Rows("1:2").Insert()

# After the insert, the former row1-4 content now lives at rows 3-6:
#   row3: A3 = "prévu"            , C3 = red swatch   (style 3), H3 = "Color Rabbit prediction "
#   row4: A4 = "Théo et Héloïse"  , C4 = gradient swatch (style5), G4 = "Planning du projet, diagramme de Grantt"
#   row5: A5 = "Théo Lemaire"     , C5 = blue swatch  (style 1)
#   row6: A6 = "Héloïse Faure"   , C6 = green swatch (style 2), J6 = "Oral 1" (rich text)
#   row7: empty
#
# We now need to reshuffle this into the new header block (rows 1-5):
#   row1: A1 = "séances prévues au début du projet" | D1 = red swatch
#   row2: A2 = "ce qui est réellement fait"          | H2 = "Color Rabbit prediction "
#   row3: C3 = "Théo et Hélo"  | D3 = gradient swatch | G3 = "Planning du projet, diagramme de Grantt"
#   row4: C4 = "Théo"          | D4 = blue swatch
#   row5: C5 = "Héloïse"       | D5 = green swatch    | J5 = "Oral 1" (rich text)

# ---------------------------------------------------------------------------
# 2) Relocate the plain-text cells (simple Cut preserves any rich runs).
#    Cut() only relocates the value, so the vacated source cell is
#    explicitly cleared afterwards (real Excel blanks it completely too).
# ---------------------------------------------------------------------------
$ws.Range("H3").Cut($ws.Range("H2"))
$ws.Range("H3").Clear()
$ws.Range("G4").Cut($ws.Range("G3"))
$ws.Range("G4").Clear()
$ws.Range("J6").Cut($ws.Range("J5"))
$ws.Range("J6").Clear()

# ---------------------------------------------------------------------------
# 3) Relocate the colour-swatch cells (Cut preserves their fill/font style).
# ---------------------------------------------------------------------------
$ws.Range("C3").Cut($ws.Range("D1"))
$ws.Range("C3").Clear()
$ws.Range("C4").Cut($ws.Range("D3"))
$ws.Range("C4").Clear()
$ws.Range("C5").Cut($ws.Range("D4"))
$ws.Range("C5").Clear()
$ws.Range("C6").Cut($ws.Range("D5"))
$ws.Range("C6").Clear()

# ---------------------------------------------------------------------------
# 4) Clear the leftover legend labels and write the new header text. The
#    assignment order below (bottom-up) matches the order the new strings
#    were originally typed in, so they land in the shared-string table in
#    the same sequence.
# ---------------------------------------------------------------------------
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("A6").ClearContents()

$ws.Range("C5").Value = "Héloïse"
$ws.Range("C4").Value = "Théo"
$ws.Range("C3").Value = "Théo et Hélo"
$ws.Range("A1").Value = "séances prévues au début du projet"
$ws.Range("A2").Value = "ce qui est réellement fait"

# ---------------------------------------------------------------------------
# 5) Draw the thin box borders — order matters, it controls the order new
#    cellXfs combinations are minted in (matches the target style table).
# ---------------------------------------------------------------------------
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("B2").Borders.LineStyle = 1
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C3").Borders.LineStyle = 1
$ws.Range("D3").Borders.LineStyle = 1
$ws.Range("C4").Borders.LineStyle = 1
$ws.Range("D4").Borders.LineStyle = 1
$ws.Range("C5").Borders.LineStyle = 1
$ws.Range("D5").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 6) Nudge the two free-floating drawing shapes back to their final spot.
# ---------------------------------------------------------------------------
$shp1 = $ws.Shapes.Item(1)
$shp1.Top = 75
$shp1.Height = 46.5

$shp2 = $ws.Shapes.Item(2)
$shp2.Top = 285

# ---------------------------------------------------------------------------
# 7) Restore the view: scrolled to A6, with F26 selected.
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A6"))
$ws.Range("F26").Select()
